# Applies the "error solve ifrs list" correction to the company_list sheet:
#   - rows 2-6 (FY2014-2018 IFRS-consolidated figures) are corrected to their
#     restated values, with a couple of now-obsolete columns cleared out
#   - rows 7-9 (the 2019E/2020E/2021E forecast rows) have all of their figures
#     cleared, keeping only the row index / ticker / period label
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: corrected figures
$ws.Range("D2").Value = 797
$ws.Range("E2").Value = -115
$ws.Range("F2").Value = -113
$ws.Range("G2").Value = -127
$ws.Range("H2").Value = -127
$ws.Range("I2").Value = -127
$ws.Range("K2").Value = 1368
$ws.Range("L2").Value = 911
$ws.Range("M2").Value = 458
$ws.Range("N2").Value = 458
$ws.Range("P2").Value = 597
$ws.Range("Q2").Value = -328
$ws.Range("R2").Value = -159
$ws.Range("S2").Value = 467
$ws.Range("T2").Value = 127
$ws.Range("U2").Value = -456
$ws.Range("V2").Value = 567
$ws.Range("W2").Value = -14.38
$ws.Range("X2").Value = -15.96
$ws.Range("Y2").Value = -24.69
$ws.Range("Z2").Value = -10.03
$ws.Range("AA2").Value = 198.94
$ws.Range("AB2").Value = -10.97
$ws.Range("AC2").Value = -168
$ws.Range("AD2").Value = -9.25
$ws.Range("AE2").Value = 603
$ws.Range("AF2").Value = 2.57
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 75915568
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3: corrected figures
$ws.Range("D3").Value = 807
$ws.Range("E3").Value = -161
$ws.Range("F3").Value = -156
$ws.Range("G3").Value = -167
$ws.Range("H3").Value = -166
$ws.Range("I3").Value = -166
$ws.Range("K3").Value = 1576
$ws.Range("L3").Value = 898
$ws.Range("M3").Value = 677
$ws.Range("N3").Value = 677
$ws.Range("P3").Value = 877
$ws.Range("Q3").Value = -195
$ws.Range("R3").Value = -143
$ws.Range("S3").Value = 396
$ws.Range("T3").Value = 44
$ws.Range("U3").Value = -238
$ws.Range("V3").Value = 611
$ws.Range("W3").Value = -19.91
$ws.Range("X3").Value = -20.58
$ws.Range("Y3").Value = -29.24
$ws.Range("Z3").Value = -11.28
$ws.Range("AA3").Value = 132.57
$ws.Range("AB3").Value = -15.03
$ws.Range("AC3").Value = -160
$ws.Range("AD3").Value = -11.86
$ws.Range("AE3").Value = 608
$ws.Range("AF3").Value = 3.13
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 111516966
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4: corrected figures
$ws.Range("D4").Value = 848
$ws.Range("E4").Value = -81
$ws.Range("F4").Value = -75
$ws.Range("G4").Value = -98
$ws.Range("H4").Value = -100
$ws.Range("I4").Value = -100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1489
$ws.Range("L4").Value = 861
$ws.Range("M4").Value = 627
$ws.Range("N4").Value = 584
$ws.Range("O4").Value = 43
$ws.Range("P4").Value = 881
$ws.Range("Q4").Value = -94
$ws.Range("R4").Value = 129
$ws.Range("S4").Value = -6
$ws.Range("T4").Value = 22
$ws.Range("U4").Value = -115
$ws.Range("V4").Value = 624
$ws.Range("W4").Value = -9.6
$ws.Range("X4").Value = -11.77
$ws.Range("Y4").Value = -15.82
$ws.Range("Z4").Value = -6.51
$ws.Range("AA4").Value = 137.27
$ws.Range("AB4").Value = -26.42
$ws.Range("AC4").Value = -89
$ws.Range("AD4").Value = -14.49
$ws.Range("AE4").Value = 521
$ws.Range("AF4").Value = 2.48
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 112039673

# Row 5: corrected figures
$ws.Range("D5").Value = 417
$ws.Range("E5").Value = -260
$ws.Range("F5").Value = -260
$ws.Range("G5").Value = -387
$ws.Range("H5").Value = -315
$ws.Range("I5").Value = -307
$ws.Range("J5").Value = -8
$ws.Range("K5").Value = 1084
$ws.Range("L5").Value = 632
$ws.Range("M5").Value = 452
$ws.Range("N5").Value = 419
$ws.Range("O5").Value = 33
$ws.Range("P5").Value = 599
$ws.Range("Q5").Value = -47
$ws.Range("R5").Value = -62
$ws.Range("S5").Value = 33
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = -52
$ws.Range("V5").Value = 509
$ws.Range("W5").Value = -62.48
$ws.Range("X5").Value = -75.61
$ws.Range("Y5").Value = -61.25
$ws.Range("Z5").Value = -24.48
$ws.Range("AA5").Value = 139.8
$ws.Range("AB5").Value = -16.64
$ws.Range("AC5").Value = -267
$ws.Range("AD5").Value = -1.83
$ws.Range("AE5").Value = 350
$ws.Range("AF5").Value = 1.39
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 119759667

# Row 6: corrected figures
$ws.Range("D6").Value = 366
$ws.Range("E6").Value = -163
$ws.Range("F6").Value = -163
$ws.Range("G6").Value = -244
$ws.Range("H6").Value = -244
$ws.Range("I6").Value = -234
$ws.Range("K6").Value = 1549
$ws.Range("L6").Value = 850
$ws.Range("M6").Value = 699
$ws.Range("N6").Value = 529
$ws.Range("P6").Value = 947
$ws.Range("Q6").Value = -134
$ws.Range("R6").Value = -99
$ws.Range("S6").Value = 473
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = -140
$ws.Range("V6").Value = 473
$ws.Range("W6").Value = -44.46
$ws.Range("X6").Value = -66.84999999999999
$ws.Range("Y6").Value = -49.31
$ws.Range("Z6").Value = -18.56
$ws.Range("AA6").Value = 121.62
$ws.Range("AB6").Value = -35.08
$ws.Range("AC6").Value = -189
$ws.Range("AD6").Value = -2.65
$ws.Range("AE6").Value = 279
$ws.Range("AF6").Value = 1.79
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 189444075

# Row 7: forecast figures removed (keep only index/ticker/period label)
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: forecast figures removed (keep only index/ticker/period label)
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: forecast figures removed (keep only index/ticker/period label)
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
